# "adicionado wait until e round"
# The scraper now waits (WebDriverWait) before reading the page, so the
# scrapy_datetime column gets stamped a bit later than before, and the
# average-price column is rounded instead of keeping its full precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated scrapy_datetime (column D) for every data row after the "wait until" was added.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "2022-05-20 19:31:54"
}

# Round the Media-Preco (column G) values to 4 decimal places.
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value2 = [Math]::Round($cell.Value2, 4)
}
